$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing columns (old O:P -> odtTilt/odtT moved into
# the new M:N slots below; deleting here shrinks the sheet's used range from
# A1:P6 down to A1:N6, matching the target dimension).
$ws.Columns("O:P").Delete()

# New header row (A1:N1) - added the "odtLoc" column and dropped the
# vcSt/vcUp/vcDn columns, reordering the remaining ones.
$headers = @("condN","condLabel","singlType","singlCont","jitTmax","stimT","gabSize","gabSf","gabNum","postStimBlankT","maskRR","odtLoc","odtTilt","odtT")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value2 = $headers[$c]
}

# New data rows, reflecting the equiluminance module's odtLoc column and the
# updated pop-out salience levels (singlCont) per condition.
$data = @{
    2 = @(1, "cont1_hi", "colour", 0.3, 500, 1000, 0.5, 0.2, 12, 250, 10, 0, 1.5, 17)
    3 = @(2, "cont1_lo", "colour", 0.3, 500, 1000, 0.5, 0.2, 12, 250, 10, 6, 1.5, 17)
    4 = @(3, "cont7_hi", "colour", 0.8, 500, 1000, 0.5, 0.2, 12, 250, 10, 0, 1.5, 17)
    5 = @(4, "cont7_lo", "colour", 0.8, 500, 1000, 0.5, 0.2, 12, 250, 10, 6, 1.5, 17)
    6 = @(5, "blank",    "colour", 0,   500, 1000, 0.5, 0.2, 12, 250, 10, 0, 1.5, 17)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}

# Move the active selection to match the authored state.
$ws.Range("E10").Select() | Out-Null
